$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 246 (Excel-style: shifts row 246..345 down to 247..346,
# and naturally grows the sheet's used range / dimension to A1:R346).
$ws.Range("A246").EntireRow.Insert()

# Populate the newly inserted row 246 with the new weekly record.
$ws.Range("A246").Value = 1
$ws.Range("B246").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C246").Value = "Arica y Parinacota"
$ws.Range("D246").Value = 44875
$ws.Range("E246").Value = 15
$ws.Range("F246").Value = 100114013
$ws.Range("G246").Value = "Zanahoria"
$ws.Range("H246").Value = "Sin especificar"
$ws.Range("I246").Value = "Primera"
$ws.Range("J246").Value = 90
$ws.Range("K246").Value = 24000
$ws.Range("L246").Value = 25000
$ws.Range("M246").Value = 24500
$ws.Range("N246").Value = "`$/saco 25 kilos"
$ws.Range("O246").Value = "Región de Arica y Parinacota"
$ws.Range("P246").Value = 980
$ws.Range("Q246").Value = 25
$ws.Range("R246").Value = "Hortaliza"
